$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows starting at row 198, shifting existing rows down
$ws.Rows.Item(198).Resize(4).EntireRow.Insert()

# Row 198
$ws.Range("A198").Value = 6
$ws.Range("B198").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C198").Value = "Metropolitana"
$ws.Range("D198").Value = 44839
$ws.Range("E198").Value = 13
$ws.Range("F198").Value = "Fruta"
$ws.Range("G198").Value = 100107
$ws.Range("H198").Value = "Otros"
$ws.Range("I198").Value = 100107002
$ws.Range("J198").Value = "Chirimoya"
$ws.Range("K198").Value = "Cultivar IV Región"
$ws.Range("L198").Value = "Especial"
$ws.Range("M198").Value = 350
$ws.Range("N198").Value = 23000
$ws.Range("O198").Value = 23000
$ws.Range("P198").Value = 23000
$ws.Range("Q198").Value = "$/bandeja 8 kilos"
$ws.Range("R198").Value = "Provincia de Limarí"
$ws.Range("S198").Value = 2875
$ws.Range("T198").Value = 8

# Row 199
$ws.Range("A199").Value = 6
$ws.Range("B199").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C199").Value = "Metropolitana"
$ws.Range("D199").Value = 44839
$ws.Range("E199").Value = 13
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100107
$ws.Range("H199").Value = "Otros"
$ws.Range("I199").Value = 100107002
$ws.Range("J199").Value = "Chirimoya"
$ws.Range("K199").Value = "Cultivar IV Región"
$ws.Range("L199").Value = "Extra (doble especial)"
$ws.Range("M199").Value = 350
$ws.Range("N199").Value = 27000
$ws.Range("O199").Value = 27000
$ws.Range("P199").Value = 27000
$ws.Range("Q199").Value = "$/bandeja 8 kilos"
$ws.Range("R199").Value = "Provincia de Limarí"
$ws.Range("S199").Value = 3375
$ws.Range("T199").Value = 8

# Row 200
$ws.Range("A200").Value = 6
$ws.Range("B200").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C200").Value = "Metropolitana"
$ws.Range("D200").Value = 44839
$ws.Range("E200").Value = 13
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100107
$ws.Range("H200").Value = "Otros"
$ws.Range("I200").Value = 100107002
$ws.Range("J200").Value = "Chirimoya"
$ws.Range("K200").Value = "Cultivar IV Región"
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 350
$ws.Range("N200").Value = 21000
$ws.Range("O200").Value = 21000
$ws.Range("P200").Value = 21000
$ws.Range("Q200").Value = "$/bandeja 8 kilos"
$ws.Range("R200").Value = "Provincia de Limarí"
$ws.Range("S200").Value = 2625
$ws.Range("T200").Value = 8

# Row 201
$ws.Range("A201").Value = 6
$ws.Range("B201").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C201").Value = "Metropolitana"
$ws.Range("D201").Value = 44839
$ws.Range("E201").Value = 13
$ws.Range("F201").Value = "Fruta"
$ws.Range("G201").Value = 100107
$ws.Range("H201").Value = "Otros"
$ws.Range("I201").Value = 100107002
$ws.Range("J201").Value = "Chirimoya"
$ws.Range("K201").Value = "Cultivar IV Región"
$ws.Range("L201").Value = "Segunda"
$ws.Range("M201").Value = 350
$ws.Range("N201").Value = 19000
$ws.Range("O201").Value = 19000
$ws.Range("P201").Value = 19000
$ws.Range("Q201").Value = "$/bandeja 8 kilos"
$ws.Range("R201").Value = "Provincia de Limarí"
$ws.Range("S201").Value = 2375
$ws.Range("T201").Value = 8
